# Generate Report for Handoff
#
# Adds two new tracked files to the localization-status report:
#   - b0089a86-f778-4424-927d-8c352a253b96.md   (a new markdown handoff, with an
#     image dependency)
#   - bdb93f05-99cf-4783-97b4-725bec8bddfc.png  (a second, independent image)
# The pre-existing row (a83ddec1-...) is superseded by a new image handoff
# 7d0d484d-7658-4242-a286-473fa110433f.png that also carries a dependency
# relationship back to the new markdown file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Row 2 already exists - update the handed-off file name / timestamp in place.
$overview.Hyperlinks.Item(1).Delete()
$overview.Hyperlinks.Add($overview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/7d0d484d-7658-4242-a286-473fa110433f.png", "", "", "7d0d484d-7658-4242-a286-473fa110433f.png")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-03-23 23:07:30"

# Row 3 - new markdown handoff.
$overview.Hyperlinks.Add($overview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/b0089a86-f778-4424-927d-8c352a253b96.md", "", "", "b0089a86-f778-4424-927d-8c352a253b96.md")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-23 23:07:30"

# Row 4 - new, unrelated image handoff.
$overview.Hyperlinks.Add($overview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/bdb93f05-99cf-4783-97b4-725bec8bddfc.png", "", "", "bdb93f05-99cf-4783-97b4-725bec8bddfc.png")
$overview.Range("B4").Value = "Ready for handoff"
$overview.Range("C4").Value = "Ready for handoff"
$overview.Range("D4").Value = "2016-03-23 23:07:30"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Hyperlinks.Item(2).Delete()
$zhcn.Hyperlinks.Item(1).Delete()

# --- Row 2: 7d0d484d-....png - now flagged as a dependency of the new .md
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/7d0d484d-7658-4242-a286-473fa110433f.png", "", "", "7d0d484d-7658-4242-a286-473fa110433f.png")
$zhcn.Range("B2").Value = ".png"
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/7396cda22bb9a402a5d706b086b711ffd5661d3f.png", "", "", "7396cda22bb9a402a5d706b086b711ffd5661d3f.png")
$zhcn.Range("E2").Value = "2016-03-23 23:07:25"
$zhcn.Range("H2").Value = "0001-01-01 00:00:00"
$zhcn.Range("J2").Value = "IsDependency"
$zhcn.Range("K2").Value = "e2e\b0089a86-f778-4424-927d-8c352a253b96.md"

# --- Row 3: b0089a86-....md - new handoff file
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/b0089a86-f778-4424-927d-8c352a253b96.md", "", "", "b0089a86-f778-4424-927d-8c352a253b96.md")
$zhcn.Range("B3").Value = ".md"
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/b0089a86-f778-4424-927d-8c352a253b96.36ec71c8d6ac1cc457ecb681fcb3a2a713450e07.zh-cn.xlf", "", "", "b0089a86-f778-4424-927d-8c352a253b96.36ec71c8d6ac1cc457ecb681fcb3a2a713450e07.zh-cn.xlf")
$zhcn.Range("E3").Value = "2016-03-23 23:07:25"
$zhcn.Range("H3").Value = "0001-01-01 00:00:00"
$zhcn.Range("J3").Value = "Include"

# --- Row 4: bdb93f05-....png - new, independent image handoff
$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/bdb93f05-99cf-4783-97b4-725bec8bddfc.png", "", "", "bdb93f05-99cf-4783-97b4-725bec8bddfc.png")
$zhcn.Range("B4").Value = ".png"
$zhcn.Range("C4").Value = "Ready for handoff"
$zhcn.Hyperlinks.Add($zhcn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/2eb685d2ac89f52383abae19e7d2a8a095672064.png", "", "", "2eb685d2ac89f52383abae19e7d2a8a095672064.png")
$zhcn.Range("E4").Value = "2016-03-23 23:07:25"
$zhcn.Range("H4").Value = "0001-01-01 00:00:00"
$zhcn.Range("J4").Value = "IsDependency"
$zhcn.Range("K4").Value = "e2e\b0089a86-f778-4424-927d-8c352a253b96.md"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Hyperlinks.Item(2).Delete()
$dede.Hyperlinks.Item(1).Delete()

# --- Row 2: 7d0d484d-....png - now flagged as a dependency of the new .md
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/7d0d484d-7658-4242-a286-473fa110433f.png", "", "", "7d0d484d-7658-4242-a286-473fa110433f.png")
$dede.Range("B2").Value = ".png"
$dede.Range("C2").Value = "Ready for handoff"
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/7396cda22bb9a402a5d706b086b711ffd5661d3f.png", "", "", "7396cda22bb9a402a5d706b086b711ffd5661d3f.png")
$dede.Range("E2").Value = "2016-03-23 23:07:30"
$dede.Range("H2").Value = "0001-01-01 00:00:00"
$dede.Range("J2").Value = "IsDependency"
$dede.Range("K2").Value = "e2e\b0089a86-f778-4424-927d-8c352a253b96.md"

# --- Row 3: b0089a86-....md - new handoff file
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/b0089a86-f778-4424-927d-8c352a253b96.md", "", "", "b0089a86-f778-4424-927d-8c352a253b96.md")
$dede.Range("B3").Value = ".md"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Hyperlinks.Add($dede.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/b0089a86-f778-4424-927d-8c352a253b96.36ec71c8d6ac1cc457ecb681fcb3a2a713450e07.de-de.xlf", "", "", "b0089a86-f778-4424-927d-8c352a253b96.36ec71c8d6ac1cc457ecb681fcb3a2a713450e07.de-de.xlf")
$dede.Range("E3").Value = "2016-03-23 23:07:30"
$dede.Range("H3").Value = "0001-01-01 00:00:00"
$dede.Range("J3").Value = "Include"

# --- Row 4: bdb93f05-....png - new, independent image handoff
$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/bdb93f05-99cf-4783-97b4-725bec8bddfc.png", "", "", "bdb93f05-99cf-4783-97b4-725bec8bddfc.png")
$dede.Range("B4").Value = ".png"
$dede.Range("C4").Value = "Ready for handoff"
$dede.Hyperlinks.Add($dede.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/2eb685d2ac89f52383abae19e7d2a8a095672064.png", "", "", "2eb685d2ac89f52383abae19e7d2a8a095672064.png")
$dede.Range("E4").Value = "2016-03-23 23:07:30"
$dede.Range("H4").Value = "0001-01-01 00:00:00"
$dede.Range("J4").Value = "IsDependency"
$dede.Range("K4").Value = "e2e\b0089a86-f778-4424-927d-8c352a253b96.md"
